$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.985.39'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '2.631.35'
$ws.Range("E3").Value = '  +3.57%  '
$ws.Range("E4").Value = '  +0.09%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '517.97'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.01%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '144.65'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.48%  '
$ws.Range("E7").Value = '  -0.24%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.565'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +0.41%  '
$ws.Range("D9").Value = '2.656.17'
$ws.Range("E9").Value = '  +4.34%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '6.26'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +2.66%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.104'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +2.33%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.336'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("E13").Value = '  -1.73%  '
$ws.Range("D14").Value = '3.101.31'
$ws.Range("E14").Value = '  +3.86%  '
$ws.Range("D15").Value = '58.945.44'
$ws.Range("E15").Value = '  +0.63%  '
$ws.Range("E16").Value = '  +0.96%  '
$ws.Range("E17").Value = '  +1.32%  '
$ws.Range("D18").Value = '2.653.88'
$ws.Range("E18").Value = '  +4.30%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '349.06'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +3.76%  '
$ws.Range("E20").Value = '  -0.20%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '10.34'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +2.40%  '
$ws.Range("E22").Value = '  +3.90%  '
$ws.Range("E23").Value = '  -0.10%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '61.67'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +1.88%  '
$ws.Range("E25").Value = '  +2.12%  '
$ws.Range("D26").Value = '2.751.34'
$ws.Range("E26").Value = '  +3.75%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.992'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.64%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.162'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.94%  '
$ws.Range("E29").Value = '  +2.25%  '
$ws.Range("E30").Value = '  +2.43%  '
$ws.Range("E31").Value = '  -0.20%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '6.26'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +7.26%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '18.98'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.28%  '
$ws.Range("E34").Value = '  +2.52%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '149.59'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.16%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.973'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +5.80%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '4.00'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +2.57%  '
$ws.Range("E38").Value = '  +2.10%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '36.73'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.80%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.843'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +2.58%  '
$ws.Range("E41").Value = '  +4.89%  '
$ws.Range("E42").Value = '  +1.35%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '278.26'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +1.93%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.995'
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0982'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -1.27%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '19.58'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +4.76%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.0527'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -1.13%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '10.29'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.12%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0229'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.23%  '
$ws.Range("D51").Value = '1.989.74'
$ws.Range("E51").Value = '  +4.26%  '
